$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$valueBlock1 = 44265.66402454182   # rows 2-15
$valueBlock2 = 44265.6424784375    # rows 16-29
$valueBlock3 = 44265.61992099537   # rows 30-43

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $valueBlock1
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $valueBlock2
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $valueBlock3
}
